$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Define a reusable "Calibri 10, theme text color" style for the new metadata row.
$dataStyle = $wb.Styles.Add("MCHDataStyle")
$dataStyle.Font.Name = "Calibri"

# Row 2: new metadata row under the header row.
# A2 = identifier, C2 = title, D2 = (blank, styled), E2 = levelOfDescription,
# F2 = extentAndMedium, G2 = notes, H2 = (blank, styled)

$a2 = $ws.Range("A2")
$a2.Value = "MCH127"
$a2.Style = "MCHDataStyle"
$a2.Font.ThemeColor = 1

$c2 = $ws.Range("C2")
$c2.Value = "NEWSPAPER CLIPPINGS, MAKANA FOOTBALL  ASSOCIATION"
$c2.Style = "MCHDataStyle"
$c2.Font.ThemeColor = 1

$d2 = $ws.Range("D2")
$d2.Style = "MCHDataStyle"
$d2.Font.ThemeColor = 1

$e2 = $ws.Range("E2")
$e2.Value = "Series"
$e2.Style = "MCHDataStyle"
$e2.Font.ThemeColor = 1

$f2 = $ws.Range("F2")
$f2.Value = "1 Box"
$f2.Style = "MCHDataStyle"
$f2.Font.ThemeColor = 1
$f2.WrapText = $false

$g2 = $ws.Range("G2")
$g2.Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"
$g2.Style = "MCHDataStyle"
$g2.Font.ThemeColor = 1

$h2 = $ws.Range("H2")
$h2.Style = "MCHDataStyle"
$h2.Font.ThemeColor = 1
